$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (columns B:AC) need to be swapped between them.
# Column A (row index) stays untouched in each row.
$rowPairs = @(
    @(9, 10),
    @(25, 26),
    @(92, 93),
    @(99, 100)
)

$firstCol = 2   # column B
$lastCol  = 29  # column AC

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}
